$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 66,2
$arr[0,0] = -0.4429532900133231
$arr[0,1] = -1.17415603059391
$arr[1,0] = -0.4101853272517285
$arr[1,1] = -1.013385614170121
$arr[2,0] = -0.3586033861970436
$arr[2,1] = -1.365544662205192
$arr[3,0] = -0.1734839704177584
$arr[3,1] = -0.8193303804456599
$arr[4,0] = 0.03513609728279504
$arr[4,1] = 0.457776570212569
$arr[5,0] = -0.3664153385402119
$arr[5,1] = -1.015889935153623
$arr[6,0] = -0.4588279046782437
$arr[6,1] = -0.9438777771958118
$arr[7,0] = -0.5088885699101326
$arr[7,1] = -1.399463548034325
$arr[8,0] = -0.5496915104726556
$arr[8,1] = -0.6387418939734986
$arr[9,0] = -0.2013099989447401
$arr[9,1] = -0.9596695935839867
$arr[10,0] = -0.2786618345183196
$arr[10,1] = -0.7735587113650648
$arr[11,0] = 0.1403098443297875
$arr[11,1] = 0.08549604127806801
$arr[12,0] = -0.1012707993499056
$arr[12,1] = -0.08894740030580389
$arr[13,0] = -0.1312857771732901
$arr[13,1] = 0.1359633195923107
$arr[14,0] = -0.222352436999347
$arr[14,1] = -0.2291607114116887
$arr[15,0] = -0.02171540917171313
$arr[15,1] = 0.6350032689937107
$arr[16,0] = 0.01678067421590068
$arr[16,1] = 0.6108036244733044
$arr[17,0] = 0.04715014032462234
$arr[17,1] = 0.5365461995492344
$arr[18,0] = -0.1100573300957268
$arr[18,1] = 0.09968061714659676
$arr[19,0] = -0.1034733988579106
$arr[19,1] = 0.09778291010786846
$arr[20,0] = 0.05956561361441697
$arr[20,1] = 0.6541943698018557
$arr[21,0] = 0.1464452615697857
$arr[21,1] = 0.5997338700898747
$arr[22,0] = 0.4926552825740285
$arr[22,1] = 1.48604415655835
$arr[23,0] = 0.1871431058305969
$arr[23,1] = 0.8648687431248027
$arr[24,0] = 0.1983401619942222
$arr[24,1] = 0.8108062811009547
$arr[25,0] = 0.1651251850079683
$arr[25,1] = 0.8084534768575925
$arr[26,0] = 0.2526248107964501
$arr[26,1] = 0.9793326399019312
$arr[27,0] = 0.4967340905935939
$arr[27,1] = 1.775261013530262
$arr[28,0] = 0.2009742246634456
$arr[28,1] = 0.8717028225393485
$arr[29,0] = 0.1452663335356389
$arr[29,1] = 0.7667251077268459
$arr[30,0] = 0.1899424260748209
$arr[30,1] = 0.94688004082663
$arr[31,0] = 0.1544143789761281
$arr[31,1] = 0.8648494824412244
$arr[32,0] = 0.1619851865203062
$arr[32,1] = 0.5718722349013543
$arr[33,0] = 0.2298942641507961
$arr[33,1] = 0.9897824747023484
$arr[34,0] = 0.1243596249140015
$arr[34,1] = 0.7094959198679706
$arr[35,0] = 0.1529327512689476
$arr[35,1] = 0.4530093282735353
$arr[36,0] = 0.4256308219826961
$arr[36,1] = 1.789954057665833
$arr[37,0] = -0.06048092184775992
$arr[37,1] = -0.2674937104707813
$arr[38,0] = 0.2909569427832806
$arr[38,1] = 0.8760650014999555
$arr[39,0] = -0.0541572794745537
$arr[39,1] = 0.6340925617217936
$arr[40,0] = 0.2035924448618995
$arr[40,1] = 1.283976095092083
$arr[41,0] = 0.2816702385643747
$arr[41,1] = 1.223316799807898
$arr[42,0] = -0.1021984228777193
$arr[42,1] = 0.008698262803146287
$arr[43,0] = -0.1517211812434963
$arr[43,1] = 0.03848301755142655
$arr[44,0] = -0.1933619869659879
$arr[44,1] = -0.4335315117798841
$arr[45,0] = -0.1951213519254005
$arr[45,1] = -0.4159742072442774
$arr[46,0] = -0.2294340797431483
$arr[46,1] = -0.4948604189408246
$arr[47,0] = -0.2249618156876537
$arr[47,1] = -0.5452151531020232
$arr[48,0] = -0.1581293522656999
$arr[48,1] = -0.3072898500840027
$arr[49,0] = -0.2276920389548281
$arr[49,1] = -0.6322535991622147
$arr[50,0] = -0.2276920389548281
$arr[50,1] = -0.6322535991622147
$arr[51,0] = -0.2056852645270732
$arr[51,1] = -0.3886905796556319
$arr[52,0] = -0.2424004902789378
$arr[52,1] = -0.5583431023766643
$arr[53,0] = -0.1693489412115966
$arr[53,1] = -0.2922315960047137
$arr[54,0] = -0.1734219249594706
$arr[54,1] = -0.3758986257486097
$arr[55,0] = -0.2286416077196036
$arr[55,1] = -0.4749304832362778
$arr[56,0] = -0.2299583523596593
$arr[56,1] = -0.6489000043300227
$arr[57,0] = -0.2746276156935938
$arr[57,1] = -0.7666348964041575
$arr[58,0] = -0.3102413385994793
$arr[58,1] = -0.8616696468257762
$arr[59,0] = -0.2307479581955122
$arr[59,1] = -0.4255184384997502
$arr[60,0] = -0.188944682991689
$arr[60,1] = -0.1413437493002201
$arr[61,0] = -0.3922015638036002
$arr[61,1] = -1.339698840908213
$arr[62,0] = -0.2899488118154345
$arr[62,1] = -0.7190782829763883
$arr[63,0] = -0.3621342162603901
$arr[63,1] = -0.9567024111461377
$arr[64,0] = -0.1376608228969651
$arr[64,1] = -0.1950721145850945
$arr[65,0] = -0.1715464908605843
$arr[65,1] = -0.3713392293581123
$ws.Range("A2:B67").Value = $arr
